$wb = $excel.ActiveWorkbook

# Update Block Detection (Short) sheet's "Front" column (E) calibration readings
$ws = $wb.Worksheets.Item("Block Detection (Short)")
$ws.Range("E3").Value = 485
$ws.Range("E4").Value = 336
$ws.Range("E5").Value = 259
$ws.Range("E6").Value = 215
$ws.Range("E7").Value = 185
$ws.Range("E8").Value = 165
$ws.Range("E9").Value = 150

# Reflect the final view state: "Block Detection (Short)" becomes the active/selected sheet
$ws.Activate()
$ws.Range("E10").Select()
